$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "['Belgium', 'Denmark', 'France']"
$ws.Range("H3").Value = "['Spain', 'West Germany', 'Portugal']"
$ws.Range("H4").Value = "['Spain', 'West Germany', 'Italy']"
$ws.Range("H5").Value = "['Netherlands', 'Soviet Union', 'Republic of Ireland']"
$ws.Range("H6").Value = "['Sweden', 'Denmark', 'France']"
$ws.Range("H7").Value = "['Germany', 'CIS', 'Netherlands']"
$ws.Range("H9").Value = "['Netherlands', 'England']"
$ws.Range("H10").Value = "['Scotland', 'England']"
$ws.Range("H11").Value = "['Netherlands', 'England']"
$ws.Range("H17").Value = "['Germany', 'Italy']"
$ws.Range("H18").Value = "['Germany', 'Czech Republic']"
$ws.Range("H19").Value = "['Germany', 'Italy']"
$ws.Range("H20").Value = "['Portugal', 'Croatia']"
$ws.Range("H25").Value = "['Belgium', 'Italy']"
$ws.Range("H26").Value = "['Turkey', 'Italy']"
$ws.Range("H27").Value = "['Spain', 'FR Yugoslavia']"
$ws.Range("H29").Value = "['Spain', 'FR Yugoslavia']"
$ws.Range("H31").Value = "['Spain', 'FR Yugoslavia']"
$ws.Range("H33").Value = "['Spain', 'FR Yugoslavia']"
$ws.Range("H35").Value = "['Spain', 'Greece']"
$ws.Range("H38").Value = "['France', 'England']"
$ws.Range("H40").Value = "['France', 'England']"
$ws.Range("H41").Value = "['Sweden', 'Denmark']"
$ws.Range("H44").Value = "['Portugal', 'Czech Republic']"
$ws.Range("H46").Value = "['Germany', 'Croatia']"
$ws.Range("H48").Value = "['Netherlands', 'Italy']"
$ws.Range("H49").Value = "['Sweden', 'Spain']"
$ws.Range("H53").Value = "['Czech Republic', 'Russia']"
$ws.Range("H54").Value = "['Germany', 'Denmark']"
$ws.Range("H57").Value = "['Spain', 'Italy']"
$ws.Range("H58").Value = "['France', 'England']"
$ws.Range("H59").Value = "['France', 'Switzerland']"
$ws.Range("H60").Value = "['Slovakia', 'England']"
$ws.Range("H61").Value = "['Wales', 'England']"
$ws.Range("H64").Value = "['Belgium', 'Italy']"
$ws.Range("H66").Value = "['Portugal', 'Hungary']"
$ws.Range("H68").Value = "['Wales', 'Italy']"
$ws.Range("H69").Value = "['Belgium', 'Finland']"
$ws.Range("H70").Value = "['Belgium', 'Denmark']"
$ws.Range("H75").Value = "['Sweden', 'Slovakia']"
$ws.Range("H76").Value = "['Sweden', 'Spain']"
$ws.Range("H81").Value = "['Germany', 'Switzerland']"
$ws.Range("H82").Value = "['Spain', 'Italy']"
$ws.Range("H84").Value = "['Spain', 'Italy']"
$ws.Range("H85").Value = "['Denmark', 'England']"
$ws.Range("H87").Value = "['Denmark', 'England']"
$ws.Range("H95").Value = "['Belgium', 'Romania']"
$ws.Range("H96").Value = "['Belgium', 'Slovakia']"
$ws.Range("H97").Value = "['Belgium', 'Romania']"
